# retraining the models for Elnet and 3D Steel
# Updates the "Forecasted Consumption (MW)" values in column A and shifts
# the "Timestamp" values in column B forward by 28 days (rows 2-97).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    5160,5120,5080,5050,5010,4980,4970,4960,4960,4960,
    4960,4970,4980,4990,5010,5040,5070,5130,5210,5310,
    5420,5540,5660,5800,5960,6070,6150,6220,6240,6230,
    6200,6130,5990,5860,5740,5620,5500,5380,5260,5170,
    5060,4980,4920,4870,4820,4780,4760,4750,4750,4750,
    4750,4760,4790,4830,4850,4900,4960,5020,5080,5150,
    5240,5340,5450,5560,5690,5820,5950,6080,6220,6340,
    6450,6570,6680,6800,6910,7010,7130,7170,7170,7130,
    7000,6850,6730,6600,6410,6260,6080,5930,5800,5650,
    5530,5420,5380,5340,5290,5220
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
    $ws.Cells.Item($row, 2).Value = $ws.Cells.Item($row, 2).Value2 + 28
}
